$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.1379634141922
$ws.Range("B1").Value = 1.844222068786621
$ws.Range("D1").Value = 2.316503524780273
$ws.Range("E1").Value = 1.111242055892944
